$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 2 new rows for "octagon" / "decagon" categories (before old row 22) ---
$ws.Rows("22:23").Insert()

$ws.Range("A22").Value = "octagon"
$ws.Range("B22").Value = "Octagon"
$ws.Range("C22").Value = 1.5

$ws.Range("A23").Value = "decagon"
$ws.Range("B23").Value = "Decagon"
$ws.Range("C23").Value = 1.5

# --- Insert 6 new rows for "prop_side_*" attributes (before what is now row 37) ---
$ws.Rows("37:42").Insert()

# Keys first (prop_side_3 .. prop_side_10), then values (Has 3 sides. .. Has 4 sides.)
# to mirror the shared-string insertion order used when this data was authored.
$ws.Range("A37").Value = "prop_side_3"
$ws.Range("A38").Value = "prop_side_4"
$ws.Range("A39").Value = "prop_side_5"
$ws.Range("A40").Value = "prop_side_6"
$ws.Range("A41").Value = "prop_side_8"
$ws.Range("A42").Value = "prop_side_10"

$ws.Range("B37").Value = "Has 3 sides."
$ws.Range("B39").Value = "Has 5 sides."
$ws.Range("B40").Value = "Has 6 sides."
$ws.Range("B41").Value = "Has 8 sides."
$ws.Range("B42").Value = "Has 10 sides."
$ws.Range("B38").Value = "Has 4 sides."

$ws.Range("C37").Value = 2
$ws.Range("C38").Value = 2
$ws.Range("C39").Value = 2
$ws.Range("C40").Value = 2
$ws.Range("C41").Value = 2
$ws.Range("C42").Value = 2

# --- Update the view state (active selection; matches the new activeCell C42) ---
$ws.Range("C42").Select()
